# The "Omschrijving" (description) column, originally column I, is moved
# so that it becomes the new column B (right after "Hoofdsoort"). All the
# columns that used to be B..H (code, theta_r, theta_s, alpha, n, lambda, Ks)
# shift one position to the right, becoming C..I.
#
# This is exactly what Excel does when you select column I, Cut it, and then
# Insert the cut cells before column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(9).Cut()
$ws.Columns.Item(2).Insert()
